$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Year_Range / Keyword / Correlation values for rows 2-21
$data = @(
    @("2008_2009", "inflation", -0.1311),
    @("2008_2009", "interest", -0.4211),
    @("2008_2009", "uncertain", -0.6423),
    @("2008_2009", "invest", -0.0548),
    @("2008_2009", "trade", -0.469),
    @("2010_2019", "uncertain", -0.1748),
    @("2010_2019", "interest", 0.4966),
    @("2010_2019", "trade", 0.0982),
    @("2010_2019", "invest", 0.4627),
    @("2010_2019", "inflation", 0.7745),
    @("2020_2021", "inflation", 0.1819),
    @("2020_2021", "interest", 0.2283),
    @("2020_2021", "invest", 0.7061),
    @("2020_2021", "trade", -0.1322),
    @("2020_2021", "uncertain", -0.1933),
    @("2022_2023", "inflation", 0.0525),
    @("2022_2023", "interest", 0.6238),
    @("2022_2023", "uncertain", -0.1727),
    @("2022_2023", "invest", -0.5367),
    @("2022_2023", "trade", 0.0973)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
